# correção nos dados e inicio da analise PNAD 2009
#
# The sheet had two problems:
#  1. Row header B2 ("sexo"/"total" sub-column) still showed the raw pandas
#     placeholder "unnamed: 1_level_1" instead of "total".
#  2. The data values in column B:F were shifted down by one row relative to
#     their row labels (column A), and two section-header rows ("situação do
#     domicílio" / "grandes regiões e unidades da federação") carried no data
#     at all. Deleting those two now-redundant rows realigns every label with
#     its correct data row and drops the trailing two rows that no longer
#     have data once everything shifts up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "total" sub-header label.
$ws.Range("B2").Value = "total"

# Row 5 ("situação do domicílio") had no data of its own - removing it moves
# the "urbana"/"rural" rows (and everything below) up so they line up with
# the values that belong to them.
$ws.Rows("5").Delete()

# After the deletion above, the former row 8 ("grandes regiões e unidades da
# federação") is now row 7. It also had no data - remove it too so "norte"
# and the remaining states/regions align with their correct values.
$ws.Rows("7").Delete()
